# ---------------------------------------------------------------
# CompStat 43rd Precinct weekly refresh: Week of 7/15-7/21/2024
# rolls forward to Week of 7/22-7/28/2024 ("Volume 31 Number 30"),
# new crime-complaint figures collected for that week.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/issue number and the covered date range ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Crime Complaints table (rows 14-30): write the refreshed values ---
$ws.Range("C14").Value = "'0"
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("N14").Value = -82.5
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("M15").Value = -19.047619047619
$ws.Range("N15").Value = -59.523809523809
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 12.5
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 26.666666666666
$ws.Range("I16").Value = 282
$ws.Range("J16").Value = 303
$ws.Range("K16").Value = -6.930693069306
$ws.Range("L16").Value = -31.219512195122
$ws.Range("M16").Value = 2.545454545454
$ws.Range("N16").Value = -75.34965034965
$ws.Range("C17").Value = 20
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = -9.090909090909
$ws.Range("F17").Value = 67
$ws.Range("G17").Value = 79
$ws.Range("H17").Value = -15.189873417721
$ws.Range("I17").Value = 423
$ws.Range("J17").Value = 497
$ws.Range("K17").Value = -14.889336016096
$ws.Range("L17").Value = -10.570824524312
$ws.Range("M17").Value = 45.360824742268
$ws.Range("N17").Value = -18.965517241379
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 19
$ws.Range("H18").Value = -42.424242424242
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 204
$ws.Range("K18").Value = -21.56862745098
$ws.Range("L18").Value = -9.604519774011
$ws.Range("M18").Value = -25.581395348837
$ws.Range("N18").Value = -84.390243902439
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 92.307692307692
$ws.Range("F19").Value = 96
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = 57.377049180327
$ws.Range("I19").Value = 568
$ws.Range("J19").Value = 538
$ws.Range("K19").Value = 5.576208178438
$ws.Range("L19").Value = -4.857621440536
$ws.Range("M19").Value = 77.5
$ws.Range("N19").Value = 34.916864608076
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 45
$ws.Range("G20").Value = 55
$ws.Range("H20").Value = -18.181818181818
$ws.Range("I20").Value = 276
$ws.Range("J20").Value = 386
$ws.Range("K20").Value = -28.497409326424
$ws.Range("L20").Value = -9.210526315789
$ws.Range("M20").Value = 102.941176470588
$ws.Range("N20").Value = -73.838862559241
$ws.Range("C21").Value = 70
$ws.Range("D21").Value = 59
$ws.Range("E21").Value = 18.64406779661
$ws.Range("F21").Value = 268
$ws.Range("G21").Value = 263
$ws.Range("H21").Value = 1.90114068441
$ws.Range("I21").Value = 1733
$ws.Range("J21").Value = 1962
$ws.Range("K21").Value = -11.671763506625
$ws.Range("L21").Value = -13.219829744616
$ws.Range("M21").Value = 37.104430379746
$ws.Range("N21").Value = -59.213932690044
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 12
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 9.090909090909
$ws.Range("M22").Value = 9.090909090909
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -44.444444444444
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = 38
$ws.Range("H23").Value = -34.210526315789
$ws.Range("I23").Value = 155
$ws.Range("J23").Value = 178
$ws.Range("K23").Value = -12.921348314606
$ws.Range("L23").Value = -23.645320197044
$ws.Range("M23").Value = 24
$ws.Range("C24").Value = 61
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 27.083333333333
$ws.Range("F24").Value = 148
$ws.Range("G24").Value = 168
$ws.Range("H24").Value = -11.904761904761
$ws.Range("I24").Value = 1040
$ws.Range("J24").Value = 1186
$ws.Range("K24").Value = -12.310286677908
$ws.Range("L24").Value = -8.851884312007
$ws.Range("M24").Value = 27.450980392156
$ws.Range("C25").Value = 28
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 55.555555555555
$ws.Range("F25").Value = 71
$ws.Range("G25").Value = 62
$ws.Range("H25").Value = 14.516129032258
$ws.Range("I25").Value = 426
$ws.Range("J25").Value = 493
$ws.Range("K25").Value = -13.590263691683
$ws.Range("L25").Value = -34.662576687116
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = 53.333333333333
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = 3.614457831325
$ws.Range("I26").Value = 683
$ws.Range("J26").Value = 662
$ws.Range("K26").Value = 3.172205438066
$ws.Range("L26").Value = 4.915514592933
$ws.Range("M26").Value = -26.082251082251
$ws.Range("C27").Value = "'0"
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = -18.918918918918
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -12.5
$ws.Range("I28").Value = 67
$ws.Range("J28").Value = 62
$ws.Range("K28").Value = 8.064516129032
$ws.Range("L28").Value = 86.111111111111
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -75
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = -42.857142857142
$ws.Range("I29").Value = 18
$ws.Range("J29").Value = 22
$ws.Range("K29").Value = -18.181818181818
$ws.Range("L29").Value = -45.454545454545
$ws.Range("M29").Value = -30.76923076923
$ws.Range("N29").Value = -78.823529411764
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = -50
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = -40
$ws.Range("I30").Value = 16
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = -5.882352941176
$ws.Range("L30").Value = -40.74074074074
$ws.Range("M30").Value = -30.434782608695

# --- A few cells flip between numeric and text representation this week
#     (e.g. a stat that was "0/0" now has data, or vice versa). Excel
#     auto-picks a generic style for the new type, so re-apply the exact
#     number format used elsewhere in the column via a formats-only paste. ---
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

